$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 158 (pushes existing rows 158..258 down to 159..259)
$ws.Rows.Item(158).Insert()

# Populate the new row 158 with the new record
$ws.Cells.Item(158, 1).Value = 4
$ws.Cells.Item(158, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(158, 3).Value = "Los Lagos"
$ws.Cells.Item(158, 4).Value = 44603
$ws.Cells.Item(158, 5).Value = 10
$ws.Cells.Item(158, 6).Value = 100112045
$ws.Cells.Item(158, 7).Value = "Zapallo"
$ws.Cells.Item(158, 8).Value = "Paine"
$ws.Cells.Item(158, 9).Value = "1a nueva(o)"
$ws.Cells.Item(158, 10).Value = 1000
$ws.Cells.Item(158, 11).Value = 500
$ws.Cells.Item(158, 12).Value = 500
$ws.Cells.Item(158, 13).Value = 500
$ws.Cells.Item(158, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(158, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(158, 16).Value = 500
$ws.Cells.Item(158, 17).Value = 1
$ws.Cells.Item(158, 18).Value = "Hortaliza"
